# Update the Icam5-Itgal sending/target cluster pairs with the new TPM-derived
# NATMI statistics. The table is expanded from the original 3x3 (minus self pairs)
# combinations to the full 3 sending clusters (ECs, FAPs, Resolving-Mac) x 4 target
# clusters (ECs, FAPs, MuSCs, Resolving-Mac) cross join, 12 rows total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Icam5","Itgal","ECs",3,1,1.523319666666667,4.569959,0.3270148067479001,0.3270148067479002,1,0.3333333333333333,0.1050723333333333,0.315217,0.006764366142259137,0.006764366142259137,0.1600587517892222,1.440528766103,0.00221204788678291,0.002212047886782911),
    @("ECs","Icam5","Itgal","FAPs",3,1,1.523319666666667,4.569959,0.3270148067479001,0.3270148067479002,3,1,0.9976189999999999,2.992857,0.06422490081252995,0.06422490081252995,1.519692642540333,13.677233782863,0.02100249352761253,0.02100249352761254),
    @("ECs","Icam5","Itgal","MuSCs",3,1,1.523319666666667,4.569959,0.3270148067479001,0.3270148067479002,3,1,0.173444,0.520332,0.01116600996625811,0.01116600996625811,0.2642106562653334,2.377895906388,0.003651450591261021,0.003651450591261022),
    @("ECs","Icam5","Itgal","Resolving-Mac",3,1,1.523319666666667,4.569959,0.3270148067479001,0.3270148067479002,3,1,14.25707666666667,42.77123,0.9178447230789528,0.9178447230789528,21.71808527550778,195.46276747957,0.3001488147422436,0.3001488147422437),
    @("FAPs","Icam5","Itgal","ECs",3,1,2.580934333333333,7.742803,0.5540555673983205,0.5540555673983206,1,0.3333333333333333,0.1050723333333333,0.315217,0.006764366142259137,0.006764366142259137,0.2711847925834445,2.440663133251,0.003747834721039374,0.003747834721039375),
    @("FAPs","Icam5","Itgal","FAPs",3,1,2.580934333333333,7.742803,0.5540555673983205,0.5540555673983206,3,1,0.9976189999999999,2.992857,0.06422490081252995,0.06422490081252995,2.574789128685667,23.173102158171,0.03558416386078714,0.03558416386078714),
    @("FAPs","Icam5","Itgal","MuSCs",3,1,2.580934333333333,7.742803,0.5540555673983205,0.5540555673983206,3,1,0.173444,0.520332,0.01116600996625811,0.01116600996625811,0.4476475745106667,4.028828170596,0.006186589987430436,0.006186589987430438),
    @("FAPs","Icam5","Itgal","Resolving-Mac",3,1,2.580934333333333,7.742803,0.5540555673983205,0.5540555673983206,3,1,14.25707666666667,42.77123,0.9178447230789528,0.9178447230789528,36.79657866196556,331.16920795769,0.5085369788290636,0.5085369788290637),
    @("Resolving-Mac","Icam5","Itgal","ECs",3,1,0.554005,1.662015,0.1189296258537793,0.1189296258537793,1,0.3333333333333333,0.1050723333333333,0.315217,0.006764366142259137,0.006764366142259137,0.05821059802833333,0.523895382255,0.0008044835344368514,0.0008044835344368514),
    @("Resolving-Mac","Icam5","Itgal","FAPs",3,1,0.554005,1.662015,0.1189296258537793,0.1189296258537793,3,1,0.9976189999999999,2.992857,0.06422490081252995,0.06422490081252995,0.5526859140949999,4.974173226854999,0.00763824342413027,0.00763824342413027),
    @("Resolving-Mac","Icam5","Itgal","MuSCs",3,1,0.554005,1.662015,0.1189296258537793,0.1189296258537793,3,1,0.173444,0.520332,0.01116600996625811,0.01116600996625811,0.09608884322,0.86479958898,0.001327969387566647,0.001327969387566647),
    @("Resolving-Mac","Icam5","Itgal","Resolving-Mac",3,1,0.554005,1.662015,0.1189296258537793,0.1189296258537793,3,1,14.25707666666667,42.77123,0.9178447230789528,0.9178447230789528,7.898491758716667,71.08642582844999,0.1091589295076455,0.1091589295076455)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
